$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is added for "Zapallo" / "Camote" dated 2021-11-30
# (serial 44530). It becomes the new row 5, pushing all the existing data
# rows (previously 5-20) down by one (to 6-21).
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 44530
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 100112045
$ws.Range("G5").Value = "Zapallo"
$ws.Range("H5").Value = "Camote"
$ws.Range("I5").Value = "2a nueva(o)"
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 480
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 490
$ws.Range("N5").Value = "$/kilo (volumen en unidades)"
$ws.Range("O5").Value = "Perú"
$ws.Range("P5").Value = 490
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
